$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "31.03.2026"
$ws.Range("B3").Value = "10:00"
$ws.Range("C3").Value = "55NM123"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "11:30"
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = "13:45"
$ws.Range("H3").Value = 14
$ws.Range("I3").Value = "14:30"
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = "Nisa Karaman"
$ws.Range("L3").Value = 9
$ws.Range("M3").Value = 10
